$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Wrong count total changed
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): Right/Wrong totals and the Max summary string changed
$ws.Range("B12").Value = 68
$ws.Range("C12").Value = -10
$ws.Range("E12").Value = "58 / 112"
